$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.329.11'
$ws.Range("E2").Value = '  +0.69%  '
$ws.Range("D3").Value = '3.822.19'
$ws.Range("E3").Value = '  -0.67%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '708.64'
$ws.Range("E5").Value = '  +1.43%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.32'
$ws.Range("E6").Value = '  -0.08%  '
$ws.Range("D7").Value = '3.821.87'
$ws.Range("E7").Value = '  -0.63%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.524'
$ws.Range("E9").Value = '  +0.02%  '
$ws.Range("E10").Value = '  -0.05%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.66'
$ws.Range("E11").Value = '  +6.08%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.462'
$ws.Range("E12").Value = '  +0.84%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000254'
$ws.Range("E13").Value = '  -1.16%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.05'
$ws.Range("E14").Value = '  -0.35%  '
$ws.Range("D15").Value = '4.467.11'
$ws.Range("E15").Value = '  -0.73%  '
$ws.Range("D16").Value = '3.821.40'
$ws.Range("E16").Value = '  -0.65%  '
$ws.Range("D17").Value = '71.316.76'
$ws.Range("E17").Value = '  +0.56%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.18'
$ws.Range("E18").Value = '  +0.12%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.52'
$ws.Range("E19").Value = '  +0.79%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '517.13'
$ws.Range("E21").Value = '  +4.49%  '
$ws.Range("E22").Value = '  +0.57%  '
$ws.Range("E23").Value = '  +0.92%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.61'
$ws.Range("E24").Value = '  -0.15%  '
$ws.Range("E25").Value = '  -2.00%  '
$ws.Range("D26").Value = '3.972.26'
$ws.Range("E26").Value = '  -0.81%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.04'
$ws.Range("E27").Value = '  -1.44%  '
$ws.Range("E28").Value = '  -1.53%  '
$ws.Range("E29").Value = '  +0.27%  '
$ws.Range("E30").Value = '  -3.40%  '
$ws.Range("E31").Value = '  -3.08%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.25'
$ws.Range("E32").Value = '  -0.33%  '
$ws.Range("B33").Value = 'NEARProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.41'
$ws.Range("E33").Value = '  -1.39%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '29.19'
$ws.Range("E34").Value = '  -0.95%  '
$ws.Range("E35").Value = '  -5.13%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.20'
$ws.Range("E36").Value = '  +0.17%  '
$ws.Range("D37").Value = '3.792.48'
$ws.Range("E37").Value = '  -0.39%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.998'
$ws.Range("E38").Value = '  -0.22%  '
$ws.Range("E39").Value = '  -1.47%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.37'
$ws.Range("E40").Value = '  -0.60%  '
$ws.Range("E41").Value = '  -0.46%  '
$ws.Range("E42").Value = '  -1.81%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.29'
$ws.Range("E43").Value = '  -2.20%  '
$ws.Range("E45").Value = '  +0.01%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '168.16'
$ws.Range("E46").Value = '  +2.55%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.000313'
$ws.Range("E47").Value = '  +0.73%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '49.41'
$ws.Range("E48").Value = '  +0.88%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '423.75'
$ws.Range("E49").Value = '  +3.36%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.65'
$ws.Range("E50").Value = '  +0.34%  '
$ws.Range("B51").Value = 'TheGraph'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.296'
$ws.Range("E51").Value = '  -1.33%  '